# Updates the cryptos list on the sheet with freshly scraped market data:
# new Price (column D) and Volume(1h) (column E) figures, and a
# position swap between the Kaspa and TheGraph rows (41/42) together
# with their respective data.
#
# The sheet stores Price/Volume as plain text (inlineStr) rather than
# numbers, so that values like "66.325.74" (a dotted/grouped price) or
# "0.740" (with a significant trailing zero) keep their exact original
# formatting instead of being reinterpreted by Excel. A handful of the
# new Price values parse as plain numbers (e.g. "595.62"); for those
# cells we briefly mark the cell as Text (NumberFormat "@") before
# assigning the value so Excel stores the literal digits instead of
# silently converting them to a number, then restore the cell style so
# no visual/style change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '66.325.74'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = '3.192.19'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '595.62'
$ws.Range("E5").Value = '  +3.30%  '
Set-TextValue $ws.Range("D6") '154.22'
$ws.Range("E6").Value = '  +3.53%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.190.36'
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("E9").Value = '  +4.16%  '
$ws.Range("E10").Value = '  +1.09%  '
Set-TextValue $ws.Range("D11") '5.99'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("E12").Value = '  +3.62%  '
Set-TextValue $ws.Range("D13") '0.0000268'
$ws.Range("E13").Value = '  +3.25%  '
Set-TextValue $ws.Range("D14") '39.08'
$ws.Range("E14").Value = '  +5.44%  '
$ws.Range("D15").Value = '3.715.56'
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").Value = '66.336.41'
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("E17").Value = '  +4.88%  '
$ws.Range("D18").Value = '3.192.08'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("E19").Value = '  +0.85%  '
Set-TextValue $ws.Range("D20") '513.88'
$ws.Range("E20").Value = '  +2.31%  '
Set-TextValue $ws.Range("D21") '15.38'
$ws.Range("E21").Value = '  +4.03%  '
Set-TextValue $ws.Range("D22") '0.740'
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("E23").Value = '  +4.89%  '
Set-TextValue $ws.Range("D24") '15.05'
$ws.Range("E24").Value = '  -1.11%  '
Set-TextValue $ws.Range("D25") '85.74'
$ws.Range("E25").Value = '  +2.34%  '
Set-TextValue $ws.Range("D26") '0.997'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +5.24%  '
Set-TextValue $ws.Range("D28") '3.01'
$ws.Range("E28").Value = '  +4.16%  '
Set-TextValue $ws.Range("D29") '2.35'
$ws.Range("E29").Value = '  +8.17%  '
Set-TextValue $ws.Range("D30") '7.15'
$ws.Range("E30").Value = '  +16.04%  '
$ws.Range("E31").Value = '  +3.06%  '
Set-TextValue $ws.Range("D32") '28.31'
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("E35").Value = '  +1.17%  '
Set-TextValue $ws.Range("D36") '509.61'
$ws.Range("E36").Value = '  +7.59%  '
Set-TextValue $ws.Range("D37") '54.84'
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  +2.33%  '
Set-TextValue $ws.Range("D40") '8.89'
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D41") '0.305'
$ws.Range("E41").Value = '  +8.54%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D42") '0.123'
$ws.Range("E42").Value = '  +7.25%  '
Set-TextValue $ws.Range("D43") '2.88'
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D44").Value = '0.0₃0672'
$ws.Range("E44").Value = '  +15.37%  '
Set-TextValue $ws.Range("D45") '2.46'
$ws.Range("E45").Value = '  +2.28%  '
$ws.Range("D46").Value = '2.917.20'
$ws.Range("E46").Value = '  -3.02%  '
Set-TextValue $ws.Range("D47") '28.73'
$ws.Range("E47").Value = '  +2.17%  '
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("E49").Value = '  +0.02%  '
Set-TextValue $ws.Range("D50") '2.34'
$ws.Range("E50").Value = '  +5.24%  '
Set-TextValue $ws.Range("D51") '2.64'
$ws.Range("E51").Value = '  +10.53%  '
